# This script applies a weekly update to the "Fruta, Vega Modelo de Temuco - Maracuyá"
# data sheet: two new rows of data (dated 2022-08-02 / serial 44775) are inserted at
# the top of the data block (row 25), pushing all the existing rows from 25..52 down
# to 27..54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 25, shifting the existing
# data (previously rows 25-52) down to rows 27-54.
$insertRange = $ws.Range("A25:T26")
$insertRange.EntireRow.Insert()

# Columns A, B, C, E, F, G, H, I, J, K are constant for every data row in this sheet.
$constA = 10
$constB = "Vega Modelo de Temuco"
$constC = "La Araucanía"
$constE = 9
$constF = "Fruta"
$constG = 100108
$constH = "Tropicales y subtropicales"
$constI = 100108003
$constJ = "Maracuyá"
$constK = "Sin especificar"

# --- New row 25 ---
$ws.Cells.Item(25, 1).Value = $constA
$ws.Cells.Item(25, 2).Value = $constB
$ws.Cells.Item(25, 3).Value = $constC
$ws.Cells.Item(25, 4).Value = 44775
$ws.Cells.Item(25, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(25, 5).Value = $constE
$ws.Cells.Item(25, 6).Value = $constF
$ws.Cells.Item(25, 7).Value = $constG
$ws.Cells.Item(25, 8).Value = $constH
$ws.Cells.Item(25, 9).Value = $constI
$ws.Cells.Item(25, 10).Value = $constJ
$ws.Cells.Item(25, 11).Value = $constK
$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 20
$ws.Cells.Item(25, 14).Value = 35000
$ws.Cells.Item(25, 15).Value = 35000
$ws.Cells.Item(25, 16).Value = 35000
$ws.Cells.Item(25, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(25, 18).Value = "Perú"
$ws.Cells.Item(25, 19).Value = 1944
$ws.Cells.Item(25, 20).Value = 18

# --- New row 26 ---
$ws.Cells.Item(26, 1).Value = $constA
$ws.Cells.Item(26, 2).Value = $constB
$ws.Cells.Item(26, 3).Value = $constC
$ws.Cells.Item(26, 4).Value = 44775
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = $constE
$ws.Cells.Item(26, 6).Value = $constF
$ws.Cells.Item(26, 7).Value = $constG
$ws.Cells.Item(26, 8).Value = $constH
$ws.Cells.Item(26, 9).Value = $constI
$ws.Cells.Item(26, 10).Value = $constJ
$ws.Cells.Item(26, 11).Value = $constK
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 40
$ws.Cells.Item(26, 14).Value = 30000
$ws.Cells.Item(26, 15).Value = 30000
$ws.Cells.Item(26, 16).Value = 30000
$ws.Cells.Item(26, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(26, 18).Value = "Perú"
$ws.Cells.Item(26, 19).Value = 1667
$ws.Cells.Item(26, 20).Value = 18
